$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 10098.962
$ws.Range("J17").Value = 10292.92
$ws.Range("L17").Value = 30878.76
$ws.Range("N17").Value = -31214.76
$ws.Range("H28").Value = 112279.78
$ws.Range("I28").Value = 134369.2
$ws.Range("J28").Value = 1832.6666
$ws.Range("K28").Value = 134369.2
$ws.Range("L28").Value = 1832.6666
$ws.Range("M28").Value = -133884.2
$ws.Range("N28").Value = -2802.6666
$ws.Range("H39").Value = 500
$ws.Range("I39").Value = 500
$ws.Range("K39").Value = 1500
$ws.Range("M39").Value = -1204
$ws.Range("H41").Value = 2168.818
$ws.Range("I41").Value = 1659
$ws.Range("J41").Value = 3061
$ws.Range("K41").Value = 1659
$ws.Range("L41").Value = 3061
$ws.Range("M41").Value = -1219
$ws.Range("N41").Value = -3941
$ws.Range("H43").Value = 2666.6667
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 2500
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 2500
$ws.Range("M43").Value = -2931
$ws.Range("N43").Value = -2638
$ws.Range("H49").Value = 1208
$ws.Range("I49").Value = 138.5
$ws.Range("K49").Value = 415.5
$ws.Range("M49").Value = -279.5
$ws.Range("H76").Value = 4433.3335
$ws.Range("I76").Value = 4400
$ws.Range("J76").Value = 4500
$ws.Range("K76").Value = 4400
$ws.Range("L76").Value = 4500
$ws.Range("M76").Value = -4085
$ws.Range("N76").Value = -5130
$ws.Range("H79").Value = 4433.3335
$ws.Range("I79").Value = 4400
$ws.Range("J79").Value = 4500
$ws.Range("K79").Value = 4400
$ws.Range("L79").Value = 4500
$ws.Range("M79").Value = -3308
$ws.Range("N79").Value = -6684
$ws.Range("H80").Value = 71437510
$ws.Range("J80").Value = 10364.083
$ws.Range("L80").Value = 31092.249
$ws.Range("N80").Value = -33088.249
$ws.Range("H83").Value = 71437510
$ws.Range("J83").Value = 10364.083
$ws.Range("L83").Value = 93276.747
$ws.Range("N83").Value = -103260.747
$ws.Range("H86").Value = 375173340
$ws.Range("I86").Value = 333346660
$ws.Range("J86").Value = 417000000
$ws.Range("K86").Value = 333346660
$ws.Range("L86").Value = 417000000
$ws.Range("M86").Value = -333345537
$ws.Range("N86").Value = -417002246
$ws.Range("H88").Value = 2863
$ws.Range("I88").Value = 2897.6667
$ws.Range("J88").Value = 2850
$ws.Range("K88").Value = 2897.6667
$ws.Range("L88").Value = 2850
$ws.Range("M88").Value = -2491.6667
$ws.Range("N88").Value = -3662
$ws.Range("H89").Value = 375173340
$ws.Range("I89").Value = 333346660
$ws.Range("J89").Value = 417000000
$ws.Range("K89").Value = 1666733300
$ws.Range("L89").Value = 2085000000
$ws.Range("M89").Value = -1666727684
$ws.Range("N89").Value = -2085011232
$ws.Range("H91").Value = 2863
$ws.Range("I91").Value = 2897.6667
$ws.Range("J91").Value = 2850
$ws.Range("K91").Value = 2897.6667
$ws.Range("L91").Value = 2850
$ws.Range("M91").Value = -1493.6667
$ws.Range("N91").Value = -5658
$ws.Range("H92").Value = 26317230
$ws.Range("I92").Value = 29413204
$ws.Range("K92").Value = 29413204
$ws.Range("M92").Value = -29411956
$ws.Range("H96").Value = 1438.6111
$ws.Range("I96").Value = 1163.4375
$ws.Range("K96").Value = 3490.3125
$ws.Range("M96").Value = -2117.3125
$ws.Range("H98").Value = 2700.625
$ws.Range("I98").Value = 2372.1428
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 2372.1428
$ws.Range("L98").Value = 5000
$ws.Range("M98").Value = -874.1428000000001
$ws.Range("N98").Value = -7996
$ws.Range("H100").Value = 3319.1538
$ws.Range("I100").Value = 1349.75
$ws.Range("K100").Value = 1349.75
$ws.Range("M100").Value = -808.75
$ws.Range("H106").Value = 9527278
$ws.Range("I106").Value = 9527278
$ws.Range("K106").Value = 9527278
$ws.Range("M106").Value = -9526647
$ws.Range("H107").Value = 1264.5
$ws.Range("I107").Value = 1268.9796
$ws.Range("K107").Value = 1268.9796
$ws.Range("M107").Value = 651.0204000000001
$ws.Range("H111").Value = 3323.8333
$ws.Range("J111").Value = 2997.6667
$ws.Range("L111").Value = 8993.000100000001
$ws.Range("N111").Value = -15127.0001
$ws.Range("H116").Value = 83366310
$ws.Range("J116").Value = 7952.6665
$ws.Range("L116").Value = 7952.6665
$ws.Range("N116").Value = -14836.6665
$ws.Range("H122").Value = 2700.625
$ws.Range("I122").Value = 2372.1428
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 7116.428400000001
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4666.428400000001
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 1235.1428
$ws.Range("I132").Value = 1222.8064
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 3668.4192
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1138.4192
$ws.Range("N132").Value = -11060
$ws.Range("H137").Value = 39533.273
$ws.Range("I137").Value = 60977.855
$ws.Range("K137").Value = 182933.565
$ws.Range("M137").Value = -180383.565
$ws.Range("H138").Value = 3474.05
$ws.Range("J138").Value = 3937.0989
$ws.Range("L138").Value = 11811.2967
$ws.Range("N138").Value = -22091.2967
$ws.Range("H141").Value = 515
$ws.Range("I141").Value = 515
$ws.Range("K141").Value = 1545
$ws.Range("M141").Value = 3635
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 5000
$ws.Range("I22").Value = 5000
$ws.Range("K22").Value = 5000
$ws.Range("M22").Value = -4701
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H32").Value = 20388306
$ws.Range("I32").Value = 20590308
$ws.Range("J32").Value = 17863266
$ws.Range("K32").Value = 20590308
$ws.Range("L32").Value = 17863266
$ws.Range("M32").Value = -20590021
$ws.Range("N32").Value = -17863840
$ws.Range("H34").Value = 35000
$ws.Range("I34").Value = 35000
$ws.Range("J34").Value = 35000
$ws.Range("K34").Value = 35000
$ws.Range("L34").Value = 35000
$ws.Range("M34").Value = -34729
$ws.Range("N34").Value = -35542
$ws.Range("H37").Value = 25034
$ws.Range("I37").Value = 25034
$ws.Range("K37").Value = 25034
$ws.Range("M37").Value = -24761
$ws.Range("H46").Value = 6775.5
$ws.Range("J46").Value = 8152
$ws.Range("L46").Value = 8152
$ws.Range("N46").Value = -8790
$ws.Range("H52").Value = 102461.75
$ws.Range("I52").Value = 80000
$ws.Range("J52").Value = 109949
$ws.Range("K52").Value = 80000
$ws.Range("L52").Value = 109949
$ws.Range("M52").Value = -79682
$ws.Range("N52").Value = -110585
$ws.Range("H55").Value = 49197
$ws.Range("I55").Value = 18000
$ws.Range("K55").Value = 18000
$ws.Range("M55").Value = -17685
$ws.Range("H61").Value = 3083.2727
$ws.Range("I61").Value = 2594.7693
$ws.Range("J61").Value = 3788.889
$ws.Range("K61").Value = 2594.7693
$ws.Range("L61").Value = 3788.889
$ws.Range("M61").Value = -2382.7693
$ws.Range("N61").Value = -4212.889
$ws.Range("H63").Value = 4986
$ws.Range("I63").Value = 4472
$ws.Range("K63").Value = 4472
$ws.Range("M63").Value = -3786
$ws.Range("H66").Value = 4986
$ws.Range("I66").Value = 4472
$ws.Range("K66").Value = 22360
$ws.Range("M66").Value = -18928
$ws.Range("H74").Value = 3102.8
$ws.Range("I74").Value = 2578.8708
$ws.Range("J74").Value = 4907.4443
$ws.Range("K74").Value = 2578.8708
$ws.Range("L74").Value = 4907.4443
$ws.Range("M74").Value = -1704.8708
$ws.Range("N74").Value = -6655.4443
$ws.Range("H77").Value = 3102.8
$ws.Range("I77").Value = 2578.8708
$ws.Range("J77").Value = 4907.4443
$ws.Range("K77").Value = 12894.354
$ws.Range("L77").Value = 24537.2215
$ws.Range("M77").Value = -8526.354000000001
$ws.Range("N77").Value = -33273.2215
$ws.Range("H102").Value = 3069.8572
$ws.Range("I102").Value = 1247.5
$ws.Range("J102").Value = 5499.6665
$ws.Range("K102").Value = 1247.5
$ws.Range("L102").Value = 5499.6665
$ws.Range("M102").Value = 374.5
$ws.Range("N102").Value = -8743.666499999999
$ws.Range("H110").Value = 2176.3076
$ws.Range("I110").Value = 1729.3
$ws.Range("J110").Value = 3666.3333
$ws.Range("K110").Value = 1729.3
$ws.Range("L110").Value = 3666.3333
$ws.Range("M110").Value = 315.7
$ws.Range("N110").Value = -7756.3333
$ws.Range("H122").Value = 3142.878
$ws.Range("I122").Value = 1638.3103
$ws.Range("J122").Value = 6778.9165
$ws.Range("K122").Value = 4914.9309
$ws.Range("L122").Value = 20336.7495
$ws.Range("M122").Value = -2464.9309
$ws.Range("N122").Value = -25236.7495
$ws.Range("H131").Value = 40000
$ws.Range("J131").Value = 40000
$ws.Range("L131").Value = 40000
$ws.Range("N131").Value = -50080
$ws.Range("H132").Value = 2411.3396
$ws.Range("I132").Value = 2093.239
$ws.Range("K132").Value = 6279.717000000001
$ws.Range("M132").Value = -3749.717000000001
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 3083.2727
$ws.Range("I136").Value = 2594.7693
$ws.Range("J136").Value = 3788.889
$ws.Range("K136").Value = 7784.3079
$ws.Range("L136").Value = 11366.667
$ws.Range("M136").Value = -5234.3079
$ws.Range("N136").Value = -16466.667
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1750.2106
$ws.Range("I20").Value = 1083.9333
$ws.Range("K20").Value = 1083.9333
$ws.Range("M20").Value = -836.9332999999999
$ws.Range("H51").Value = 59958.5
$ws.Range("J51").Value = 59958.5
$ws.Range("L51").Value = 59958.5
$ws.Range("N51").Value = -60940.5
$ws.Range("H55").Value = 86616.164
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 86616.164
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 86616.164
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -87162.164
$ws.Range("H99").Value = 3649.9473
$ws.Range("I99").Value = 2776.889
$ws.Range("J99").Value = 4435.7
$ws.Range("K99").Value = 2776.889
$ws.Range("L99").Value = 4435.7
$ws.Range("M99").Value = -1278.889
$ws.Range("N99").Value = -7431.7
$ws.Range("H105").Value = 1902.4375
$ws.Range("I105").Value = 1719.9656
$ws.Range("J105").Value = 3666.3333
$ws.Range("K105").Value = 1719.9656
$ws.Range("L105").Value = 3666.3333
$ws.Range("M105").Value = 27.03440000000001
$ws.Range("N105").Value = -7160.3333
$ws.Range("H107").Value = 2712.4
$ws.Range("I107").Value = 2285
$ws.Range("J107").Value = 2997.3333
$ws.Range("K107").Value = 2285
$ws.Range("L107").Value = 2997.3333
$ws.Range("M107").Value = -365
$ws.Range("N107").Value = -6837.3333
$ws.Range("H126").Value = 133000
$ws.Range("J126").Value = 133000
$ws.Range("L126").Value = 133000
$ws.Range("N126").Value = -142880
$ws.Range("H127").Value = 77500
$ws.Range("J127").Value = 77500
$ws.Range("L127").Value = 77500
$ws.Range("N127").Value = -87420
$ws.Range("H134").Value = 2859951
$ws.Range("I134").Value = 4203418
$ws.Range("J134").Value = 5083.625
$ws.Range("K134").Value = 12610254
$ws.Range("L134").Value = 15250.875
$ws.Range("M134").Value = -12607719
$ws.Range("N134").Value = -20320.875
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 4831428.5
$ws.Range("I6").Value = 4831428.5
$ws.Range("K6").Value = 4831428.5
$ws.Range("M6").Value = -4831315.5
$ws.Range("H16").Value = 2005.2222
$ws.Range("J16").Value = 2912
$ws.Range("L16").Value = 2912
$ws.Range("N16").Value = -3486
$ws.Range("H22").Value = 789.5
$ws.Range("I22").Value = 801.2222
$ws.Range("J22").Value = 768.4
$ws.Range("K22").Value = 801.2222
$ws.Range("L22").Value = 768.4
$ws.Range("M22").Value = -451.2222
$ws.Range("N22").Value = -1468.4
$ws.Range("H31").Value = 6229.069
$ws.Range("I31").Value = 2541.0625
$ws.Range("K31").Value = 2541.0625
$ws.Range("M31").Value = -2246.0625
$ws.Range("H34").Value = 6229.069
$ws.Range("I34").Value = 2541.0625
$ws.Range("K34").Value = 2541.0625
$ws.Range("M34").Value = -2339.0625
$ws.Range("H58").Value = 2112
$ws.Range("I58").Value = 2374.138
$ws.Range("J58").Value = 1161.75
$ws.Range("K58").Value = 2374.138
$ws.Range("L58").Value = 1161.75
$ws.Range("M58").Value = -2171.138
$ws.Range("N58").Value = -1567.75
$ws.Range("H68").Value = 50666.332
$ws.Range("J68").Value = 79999
$ws.Range("L68").Value = 79999
$ws.Range("N68").Value = -81497
$ws.Range("H71").Value = 50666.332
$ws.Range("J71").Value = 79999
$ws.Range("L71").Value = 239997
$ws.Range("N71").Value = -247485
$ws.Range("H80").Value = 100633.336
$ws.Range("J80").Value = 100633.336
$ws.Range("L80").Value = 100633.336
$ws.Range("N80").Value = -102879.336
$ws.Range("H83").Value = 100633.336
$ws.Range("J83").Value = 100633.336
$ws.Range("L83").Value = 301900.008
$ws.Range("N83").Value = -313132.008
$ws.Range("H86").Value = 42003.082
$ws.Range("I86").Value = 34620
$ws.Range("J86").Value = 43193.902
$ws.Range("K86").Value = 34620
$ws.Range("L86").Value = 43193.902
$ws.Range("M86").Value = -33497
$ws.Range("N86").Value = -45439.902
$ws.Range("H89").Value = 42003.082
$ws.Range("I89").Value = 34620
$ws.Range("J89").Value = 43193.902
$ws.Range("K89").Value = 173100
$ws.Range("L89").Value = 215969.51
$ws.Range("M89").Value = -167484
$ws.Range("N89").Value = -227201.51
$ws.Range("H94").Value = 1445
$ws.Range("I94").Value = 764
$ws.Range("J94").Value = 1626.6
$ws.Range("K94").Value = 764
$ws.Range("L94").Value = 1626.6
$ws.Range("M94").Value = -313
$ws.Range("N94").Value = -2528.6
$ws.Range("H99").Value = 3478.8
$ws.Range("I99").Value = 3473.5
$ws.Range("K99").Value = 3473.5
$ws.Range("M99").Value = -1975.5
$ws.Range("H106").Value = 50911.555
$ws.Range("I106").Value = 57999.5
$ws.Range("J106").Value = 48886.43
$ws.Range("K106").Value = 57999.5
$ws.Range("L106").Value = 48886.43
$ws.Range("M106").Value = -56737.5
$ws.Range("N106").Value = -51410.43
$ws.Range("H107").Value = 64147.75
$ws.Range("I107").Value = 112024.445
$ws.Range("J107").Value = 2592
$ws.Range("K107").Value = 112024.445
$ws.Range("L107").Value = 2592
$ws.Range("M107").Value = -110104.445
$ws.Range("N107").Value = -6432
$ws.Range("H111").Value = 48657.8
$ws.Range("J111").Value = 48657.8
$ws.Range("L111").Value = 48657.8
$ws.Range("N111").Value = -56837.8
$ws.Range("H113").Value = 2005.2222
$ws.Range("J113").Value = 2912
$ws.Range("L113").Value = 2912
$ws.Range("N113").Value = -7252
$ws.Range("H120").Value = 39999
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 16673996
$ws.Range("I122").Value = 20007196
$ws.Range("J122").Value = 7990
$ws.Range("K122").Value = 60021588
$ws.Range("L122").Value = 23970
$ws.Range("M122").Value = -60019138
$ws.Range("N122").Value = -28870
$ws.Range("H125").Value = 82652.5
$ws.Range("J125").Value = 82652.5
$ws.Range("L125").Value = 82652.5
$ws.Range("N125").Value = -87572.5
$ws.Range("H126").Value = 3478.8
$ws.Range("I126").Value = 3473.5
$ws.Range("K126").Value = 10420.5
$ws.Range("M126").Value = -7950.5
$ws.Range("H132").Value = 5368.72
$ws.Range("I132").Value = 5322.609
$ws.Range("K132").Value = 15967.827
$ws.Range("M132").Value = -13437.827
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 2286.7097
$ws.Range("I134").Value = 2119.7
$ws.Range("J134").Value = 2590.3635
$ws.Range("K134").Value = 6359.099999999999
$ws.Range("L134").Value = 7771.0905
$ws.Range("M134").Value = -3824.099999999999
$ws.Range("N134").Value = -12841.0905
$ws.Range("H136").Value = 2112
$ws.Range("I136").Value = 2374.138
$ws.Range("J136").Value = 1161.75
$ws.Range("K136").Value = 7122.414
$ws.Range("L136").Value = 3485.25
$ws.Range("M136").Value = -4572.414
$ws.Range("N136").Value = -8585.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 30.777779
$ws.Range("I2").Value = 30.333334
$ws.Range("K2").Value = 182.000004
$ws.Range("M2").Value = -69.00000399999999
$ws.Range("H3").Value = 3985.7144
$ws.Range("J3").Value = 1500
$ws.Range("L3").Value = 4500
$ws.Range("N3").Value = -4724
$ws.Range("I4").Value = 176100100
$ws.Range("J4").Value = 79948860
$ws.Range("K4").Value = 528300300
$ws.Range("L4").Value = 239846580
$ws.Range("M4").Value = -528300188
$ws.Range("N4").Value = -239846804
$ws.Range("H5").Value = 3352.8
$ws.Range("I5").Value = 1741
$ws.Range("K5").Value = 5223
$ws.Range("M5").Value = -5111
$ws.Range("H7").Value = 794.7059
$ws.Range("I7").Value = 531
$ws.Range("J7").Value = 920.8261
$ws.Range("K7").Value = 1593
$ws.Range("L7").Value = 2762.4783
$ws.Range("M7").Value = -1481
$ws.Range("N7").Value = -2986.4783
$ws.Range("H33").Value = 893.4286
$ws.Range("J33").Value = 893.4286
$ws.Range("L33").Value = 5360.571599999999
$ws.Range("N33").Value = -5926.571599999999
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H92").Value = 1163.375
$ws.Range("I92").Value = 1167
$ws.Range("J92").Value = 1161.2
$ws.Range("K92").Value = 3501
$ws.Range("L92").Value = 3483.6
$ws.Range("M92").Value = -2253
$ws.Range("N92").Value = -5979.6
$ws.Range("H107").Value = 781.931
$ws.Range("I107").Value = 1821.75
$ws.Range("J107").Value = 615.5599999999999
$ws.Range("K107").Value = 5465.25
$ws.Range("L107").Value = 1846.68
$ws.Range("M107").Value = -3545.25
$ws.Range("N107").Value = -5686.68
$ws.Range("H110").Value = 5000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 15000
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -23180
$ws.Range("H122").Value = 436382.2
$ws.Range("J122").Value = 775102
$ws.Range("L122").Value = 6975918
$ws.Range("N122").Value = -6980818
$ws.Range("H130").Value = 3343.1667
$ws.Range("I130").Value = 2514.75
$ws.Range("K130").Value = 7544.25
$ws.Range("M130").Value = -2524.25
$ws.Range("H131").Value = 1270.9474
$ws.Range("J131").Value = 1709.1
$ws.Range("L131").Value = 5127.299999999999
$ws.Range("N131").Value = -15207.3
$ws.Range("H135").Value = 3352.8
$ws.Range("I135").Value = 1741
$ws.Range("K135").Value = 15669
$ws.Range("M135").Value = -13134
$ws.Range("H136").Value = 685.3
$ws.Range("I136").Value = 685.3
$ws.Range("K136").Value = 2055.9
$ws.Range("M136").Value = 3044.1
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 3745.5557
$ws.Range("I2").Value = 39
$ws.Range("K2").Value = 39
$ws.Range("M2").Value = 74
$ws.Range("H39").Value = 33000
$ws.Range("J39").Value = 33000
$ws.Range("L39").Value = 33000
$ws.Range("N39").Value = -34064
$ws.Range("H51").Value = 90978
$ws.Range("J51").Value = 90978
$ws.Range("L51").Value = 90978
$ws.Range("N51").Value = -91996
$ws.Range("H102").Value = 3576.0625
$ws.Range("I102").Value = 3515.5715
$ws.Range("K102").Value = 3515.5715
$ws.Range("M102").Value = -1893.5715
$ws.Range("H120").Value = 79923.25
$ws.Range("J120").Value = 79923.25
$ws.Range("L120").Value = 79923.25
$ws.Range("N120").Value = -89599.25
$ws.Range("H121").Value = 89653.664
$ws.Range("J121").Value = 89653.664
$ws.Range("L121").Value = 89653.664
$ws.Range("N121").Value = -93147.664
$ws.Range("H122").Value = 3722.65
$ws.Range("I122").Value = 4597.8335
$ws.Range("J122").Value = 2409.875
$ws.Range("K122").Value = 13793.5005
$ws.Range("L122").Value = 7229.625
$ws.Range("M122").Value = -11343.5005
$ws.Range("N122").Value = -12129.625
$ws.Range("H125").Value = 85097.8
$ws.Range("J125").Value = 85097.8
$ws.Range("L125").Value = 85097.8
$ws.Range("N125").Value = -90017.8
$ws.Range("H132").Value = 3281.5588
$ws.Range("I132").Value = 2707.3215
$ws.Range("J132").Value = 5961.3335
$ws.Range("K132").Value = 8121.9645
$ws.Range("L132").Value = 17884.0005
$ws.Range("M132").Value = -5591.9645
$ws.Range("N132").Value = -22944.0005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6036.6553
$ws.Range("I7").Value = 6242.3887
$ws.Range("J7").Value = 5700
$ws.Range("K7").Value = 6242.3887
$ws.Range("L7").Value = 5700
$ws.Range("M7").Value = -6130.3887
$ws.Range("N7").Value = -5924
$ws.Range("H9").Value = 360
$ws.Range("I9").Value = 433
$ws.Range("J9").Value = 141
$ws.Range("K9").Value = 433
$ws.Range("L9").Value = 141
$ws.Range("M9").Value = -209
$ws.Range("N9").Value = -589
$ws.Range("H22").Value = 3723.1052
$ws.Range("I22").Value = 3211
$ws.Range("J22").Value = 4832.6665
$ws.Range("K22").Value = 3211
$ws.Range("L22").Value = 4832.6665
$ws.Range("M22").Value = -2916
$ws.Range("N22").Value = -5422.6665
$ws.Range("H27").Value = 3723.1052
$ws.Range("I27").Value = 3211
$ws.Range("J27").Value = 4832.6665
$ws.Range("K27").Value = 3211
$ws.Range("L27").Value = 4832.6665
$ws.Range("M27").Value = -3104
$ws.Range("N27").Value = -5046.6665
$ws.Range("H29").Value = 48495
$ws.Range("J29").Value = 48495
$ws.Range("L29").Value = 48495
$ws.Range("N29").Value = -49085
$ws.Range("H35").Value = 5033.3335
$ws.Range("I35").Value = 5033.3335
$ws.Range("K35").Value = 5033.3335
$ws.Range("M35").Value = -4697.3335
$ws.Range("H38").Value = 60040
$ws.Range("J38").Value = 60040
$ws.Range("L38").Value = 60040
$ws.Range("N38").Value = -60860
$ws.Range("H40").Value = 16671662
$ws.Range("I40").Value = 20836764
$ws.Range("K40").Value = 20836764
$ws.Range("M40").Value = -20836628
$ws.Range("H46").Value = 2165.8333
$ws.Range("I46").Value = 998.75
$ws.Range("J46").Value = 4500
$ws.Range("K46").Value = 998.75
$ws.Range("L46").Value = 4500
$ws.Range("M46").Value = -810.75
$ws.Range("N46").Value = -4876
$ws.Range("H55").Value = 517.8333
$ws.Range("J55").Value = 897
$ws.Range("L55").Value = 897
$ws.Range("N55").Value = -1243
$ws.Range("H61").Value = 1560.7142
$ws.Range("I61").Value = 1735
$ws.Range("J61").Value = 921.6667
$ws.Range("K61").Value = 1735
$ws.Range("L61").Value = 921.6667
$ws.Range("M61").Value = -1533
$ws.Range("N61").Value = -1325.6667
$ws.Range("H68").Value = 3956
$ws.Range("I68").Value = 3369.4482
$ws.Range("J68").Value = 5373.5
$ws.Range("K68").Value = 3369.4482
$ws.Range("L68").Value = 5373.5
$ws.Range("M68").Value = -2620.4482
$ws.Range("N68").Value = -6871.5
$ws.Range("H71").Value = 3956
$ws.Range("I71").Value = 3369.4482
$ws.Range("J71").Value = 5373.5
$ws.Range("K71").Value = 16847.241
$ws.Range("L71").Value = 26867.5
$ws.Range("M71").Value = -13103.241
$ws.Range("N71").Value = -34355.5
$ws.Range("H75").Value = 103323
$ws.Range("J75").Value = 103323
$ws.Range("L75").Value = 103323
$ws.Range("N75").Value = -105195
$ws.Range("H78").Value = 103323
$ws.Range("J78").Value = 103323
$ws.Range("L78").Value = 309969
$ws.Range("N78").Value = -319329
$ws.Range("H93").Value = 62501580
$ws.Range("I93").Value = 83334140
$ws.Range("K93").Value = 83334140
$ws.Range("M93").Value = -83332892
$ws.Range("H113").Value = 1560.7142
$ws.Range("I113").Value = 1735
$ws.Range("J113").Value = 921.6667
$ws.Range("K113").Value = 1735
$ws.Range("L113").Value = 921.6667
$ws.Range("M113").Value = 435
$ws.Range("N113").Value = -5261.6667
$ws.Range("H120").Value = 122999
$ws.Range("J120").Value = 122999
$ws.Range("L120").Value = 122999
$ws.Range("N120").Value = -132675
$ws.Range("H121").Value = 65974.5
$ws.Range("J121").Value = 65974.5
$ws.Range("L121").Value = 65974.5
$ws.Range("N121").Value = -69468.5
$ws.Range("H122").Value = 7643.4136
$ws.Range("I122").Value = 6985.8335
$ws.Range("J122").Value = 10799.8
$ws.Range("K122").Value = 20957.5005
$ws.Range("L122").Value = 32399.4
$ws.Range("M122").Value = -18507.5005
$ws.Range("N122").Value = -37299.39999999999
$ws.Range("H126").Value = 6036.6553
$ws.Range("I126").Value = 6242.3887
$ws.Range("J126").Value = 5700
$ws.Range("K126").Value = 18727.1661
$ws.Range("L126").Value = 17100
$ws.Range("M126").Value = -16257.1661
$ws.Range("N126").Value = -22040
$ws.Range("H132").Value = 35215
$ws.Range("I132").Value = 42569.594
$ws.Range("J132").Value = 6847.2856
$ws.Range("K132").Value = 127708.782
$ws.Range("L132").Value = 20541.8568
$ws.Range("M132").Value = -125178.782
$ws.Range("N132").Value = -25601.8568
$ws.Range("H135").Value = 53868.42
$ws.Range("J135").Value = 53868.42
$ws.Range("L135").Value = 53868.42
$ws.Range("N135").Value = -64008.42
$ws.Range("H136").Value = 2625.0667
$ws.Range("I136").Value = 2067.2
$ws.Range("K136").Value = 6201.599999999999
$ws.Range("M136").Value = -3651.599999999999
$ws.Range("H141").Value = 549999.5
$ws.Range("J141").Value = 549999.5
$ws.Range("L141").Value = 549999.5
$ws.Range("N141").Value = -560359.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 114992.336
$ws.Range("J16").Value = 114992.336
$ws.Range("L16").Value = 114992.336
$ws.Range("N16").Value = -115576.336
$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 10000
$ws.Range("K29").Value = 10000
$ws.Range("M29").Value = -9710
$ws.Range("H32").Value = 12112.75
$ws.Range("I32").Value = 12112.75
$ws.Range("K32").Value = 12112.75
$ws.Range("M32").Value = -11795.75
$ws.Range("H64").Value = 95987.2
$ws.Range("I64").Value = 89996
$ws.Range("K64").Value = 89996
$ws.Range("M64").Value = -89748
$ws.Range("H67").Value = 95987.2
$ws.Range("I67").Value = 89996
$ws.Range("K67").Value = 89996
$ws.Range("M67").Value = -89138
$ws.Range("H81").Value = 3468.3684
$ws.Range("I81").Value = 2921.4285
$ws.Range("K81").Value = 5842.857
$ws.Range("M81").Value = -4781.857
$ws.Range("H84").Value = 3468.3684
$ws.Range("I84").Value = 2921.4285
$ws.Range("K84").Value = 29214.285
$ws.Range("M84").Value = -23910.285
$ws.Range("H95").Value = 58992.5
$ws.Range("J95").Value = 58992.5
$ws.Range("L95").Value = 58992.5
$ws.Range("N95").Value = -64484.5
$ws.Range("H104").Value = 54806.332
$ws.Range("J104").Value = 54806.332
$ws.Range("L104").Value = 54806.332
$ws.Range("N104").Value = -61794.332
$ws.Range("H107").Value = 680.625
$ws.Range("I107").Value = 489.2
$ws.Range("K107").Value = 1467.6
$ws.Range("M107").Value = 452.4000000000001
$ws.Range("H118").Value = 104529.664
$ws.Range("J118").Value = 104529.664
$ws.Range("L118").Value = 104529.664
$ws.Range("N118").Value = -107843.664
$ws.Range("H119").Value = 140000
$ws.Range("J119").Value = 140000
$ws.Range("L119").Value = 140000
$ws.Range("N119").Value = -149676
$ws.Range("H120").Value = 93665.14
$ws.Range("J120").Value = 93665.14
$ws.Range("L120").Value = 93665.14
$ws.Range("N120").Value = -103341.14
$ws.Range("H121").Value = 74577.60000000001
$ws.Range("J121").Value = 74577.60000000001
$ws.Range("L121").Value = 74577.60000000001
$ws.Range("N121").Value = -78071.60000000001
$ws.Range("H122").Value = 62504460
$ws.Range("I122").Value = 200000530
$ws.Range("J122").Value = 6248.091
$ws.Range("K122").Value = 600001590
$ws.Range("L122").Value = 18744.273
$ws.Range("M122").Value = -599999140
$ws.Range("N122").Value = -23644.273
$ws.Range("H126").Value = 3787.611
$ws.Range("I126").Value = 4125.1333
$ws.Range("J126").Value = 2100
$ws.Range("K126").Value = 12375.3999
$ws.Range("L126").Value = 6300
$ws.Range("M126").Value = -9905.3999
$ws.Range("N126").Value = -11240
$ws.Range("H132").Value = 2215.5256
$ws.Range("I132").Value = 1852
$ws.Range("K132").Value = 5556
$ws.Range("M132").Value = -3026
$ws.Range("H136").Value = 42176.84
$ws.Range("I136").Value = 1509.5
$ws.Range("J136").Value = 146750
$ws.Range("K136").Value = 4528.5
$ws.Range("L136").Value = 440250
$ws.Range("M136").Value = -1978.5
$ws.Range("N136").Value = -445350
